$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change, name stays
$ws.Range("B3").Value = 0.9504773626993493
$ws.Range("C3").Value = 0.9454408164439498
$ws.Range("D3").Value = 0.9270231819066974

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.8754842061590548
$ws.Range("C4").Value = 0.8791851613980457
$ws.Range("D4").Value = 0.6432367287055977

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8966433969623976
$ws.Range("C5").Value = 0.8974724030981669
$ws.Range("D5").Value = 0.8916533234895564
